$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(18, 8).Value = 58826132
$ws.Cells.Item(18, 9).Value = 2267.6155
$ws.Cells.Item(18, 11).Value = 2267.6155
$ws.Cells.Item(18, 13).Value = -1983.6155
$ws.Cells.Item(28, 8).Value = 1317.2273
$ws.Cells.Item(28, 9).Value = 440.35294
$ws.Cells.Item(28, 10).Value = 4298.6
$ws.Cells.Item(28, 11).Value = 440.35294
$ws.Cells.Item(28, 12).Value = 4298.6
$ws.Cells.Item(28, 13).Value = 44.64706000000001
$ws.Cells.Item(28, 14).Value = -5268.6
$ws.Cells.Item(43, 8).Value = 3036.25
$ws.Cells.Item(43, 9).Value = 3900
$ws.Cells.Item(43, 10).Value = 2172.5
$ws.Cells.Item(43, 11).Value = 3900
$ws.Cells.Item(43, 12).Value = 2172.5
$ws.Cells.Item(43, 13).Value = -3831
$ws.Cells.Item(43, 14).Value = -2310.5
$ws.Cells.Item(55, 8).Value = 521.5
$ws.Cells.Item(55, 10).Value = 206.33333
$ws.Cells.Item(55, 12).Value = 206.33333
$ws.Cells.Item(55, 14).Value = -634.3333299999999
$ws.Cells.Item(62, 8).Value = 2752.5
$ws.Cells.Item(62, 9).Value = 2742.8667
$ws.Cells.Item(62, 11).Value = 2742.8667
$ws.Cells.Item(62, 13).Value = -2118.8667
$ws.Cells.Item(65, 8).Value = 2752.5
$ws.Cells.Item(65, 9).Value = 2742.8667
$ws.Cells.Item(65, 11).Value = 13714.3335
$ws.Cells.Item(65, 13).Value = -10594.3335
$ws.Cells.Item(86, 8).Value = 333334180
$ws.Cells.Item(86, 9).Value = 333334180
$ws.Cells.Item(86, 11).Value = 333334180
$ws.Cells.Item(86, 13).Value = -333333057
$ws.Cells.Item(89, 8).Value = 333334180
$ws.Cells.Item(89, 9).Value = 333334180
$ws.Cells.Item(89, 11).Value = 1666670900
$ws.Cells.Item(89, 13).Value = -1666665284
$ws.Cells.Item(92, 8).Value = 570.2727
$ws.Cells.Item(92, 9).Value = 570.2727
$ws.Cells.Item(92, 11).Value = 570.2727
$ws.Cells.Item(92, 13).Value = 677.7273
$ws.Cells.Item(96, 8).Value = 722
$ws.Cells.Item(96, 9).Value = 803.6667
$ws.Cells.Item(96, 11).Value = 2411.0001
$ws.Cells.Item(96, 13).Value = -1038.0001
$ws.Cells.Item(106, 8).Value = 1716.5555
$ws.Cells.Item(106, 9).Value = 1543.625
$ws.Cells.Item(106, 11).Value = 1543.625
$ws.Cells.Item(106, 13).Value = -912.625
$ws.Cells.Item(137, 8).Value = 8620.643
$ws.Cells.Item(137, 9).Value = 1697.7
$ws.Cells.Item(137, 10).Value = 25928
$ws.Cells.Item(137, 11).Value = 5093.1
$ws.Cells.Item(137, 12).Value = 77784
$ws.Cells.Item(137, 13).Value = -2543.1
$ws.Cells.Item(137, 14).Value = -82884

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 3649.18
$ws.Cells.Item(32, 9).Value = 3649.18
$ws.Cells.Item(32, 11).Value = 3649.18
$ws.Cells.Item(32, 13).Value = -3362.18
$ws.Cells.Item(61, 8).Value = 3985.4119
$ws.Cells.Item(61, 9).Value = 2254.75
$ws.Cells.Item(61, 10).Value = 5523.778
$ws.Cells.Item(61, 11).Value = 2254.75
$ws.Cells.Item(61, 12).Value = 5523.778
$ws.Cells.Item(61, 13).Value = -2042.75
$ws.Cells.Item(61, 14).Value = -5947.778
$ws.Cells.Item(74, 8).Value = 215405.84
$ws.Cells.Item(74, 9).Value = 309896.06
$ws.Cells.Item(74, 10).Value = 2802.875
$ws.Cells.Item(74, 11).Value = 309896.06
$ws.Cells.Item(74, 12).Value = 2802.875
$ws.Cells.Item(74, 13).Value = -309022.06
$ws.Cells.Item(74, 14).Value = -4550.875
$ws.Cells.Item(77, 8).Value = 215405.84
$ws.Cells.Item(77, 9).Value = 309896.06
$ws.Cells.Item(77, 10).Value = 2802.875
$ws.Cells.Item(77, 11).Value = 1549480.3
$ws.Cells.Item(77, 12).Value = 14014.375
$ws.Cells.Item(77, 13).Value = -1545112.3
$ws.Cells.Item(77, 14).Value = -22750.375
$ws.Cells.Item(88, 8).Value = 3851.923
$ws.Cells.Item(88, 9).Value = 1868.8572
$ws.Cells.Item(88, 11).Value = 1868.8572
$ws.Cells.Item(88, 13).Value = -1462.8572
$ws.Cells.Item(91, 8).Value = 3851.923
$ws.Cells.Item(91, 9).Value = 1868.8572
$ws.Cells.Item(91, 11).Value = 1868.8572
$ws.Cells.Item(91, 13).Value = -464.8571999999999
$ws.Cells.Item(122, 8).Value = 4979.591
$ws.Cells.Item(122, 9).Value = 4685.421
$ws.Cells.Item(122, 10).Value = 6842.6665
$ws.Cells.Item(122, 11).Value = 14056.263
$ws.Cells.Item(122, 12).Value = 20527.9995
$ws.Cells.Item(122, 13).Value = -11606.263
$ws.Cells.Item(122, 14).Value = -25427.9995
$ws.Cells.Item(135, 8).Value = 106831.75
$ws.Cells.Item(135, 10).Value = 106831.75
$ws.Cells.Item(135, 12).Value = 106831.75
$ws.Cells.Item(135, 14).Value = -116971.75
$ws.Cells.Item(136, 8).Value = 3985.4119
$ws.Cells.Item(136, 9).Value = 2254.75
$ws.Cells.Item(136, 10).Value = 5523.778
$ws.Cells.Item(136, 11).Value = 6764.25
$ws.Cells.Item(136, 12).Value = 16571.334
$ws.Cells.Item(136, 13).Value = -4214.25
$ws.Cells.Item(136, 14).Value = -21671.334

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(134, 8).Value = 2506.8
$ws.Cells.Item(134, 9).Value = 2277.077
$ws.Cells.Item(134, 11).Value = 6831.231000000001
$ws.Cells.Item(134, 13).Value = -4296.231000000001

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(16, 8).Value = 1305.8223
$ws.Cells.Item(16, 9).Value = 1214.6875
$ws.Cells.Item(16, 11).Value = 1214.6875
$ws.Cells.Item(16, 13).Value = -927.6875
$ws.Cells.Item(31, 8).Value = 6287.95
$ws.Cells.Item(31, 10).Value = 8808.125
$ws.Cells.Item(31, 12).Value = 8808.125
$ws.Cells.Item(31, 14).Value = -9398.125
$ws.Cells.Item(34, 8).Value = 6287.95
$ws.Cells.Item(34, 10).Value = 8808.125
$ws.Cells.Item(34, 12).Value = 8808.125
$ws.Cells.Item(34, 14).Value = -9212.125
$ws.Cells.Item(58, 8).Value = 2538.3125
$ws.Cells.Item(58, 9).Value = 1656.2858
$ws.Cells.Item(58, 10).Value = 4222.1816
$ws.Cells.Item(58, 11).Value = 1656.2858
$ws.Cells.Item(58, 12).Value = 4222.1816
$ws.Cells.Item(58, 13).Value = -1453.2858
$ws.Cells.Item(58, 14).Value = -4628.1816
$ws.Cells.Item(113, 8).Value = 1305.8223
$ws.Cells.Item(113, 9).Value = 1214.6875
$ws.Cells.Item(113, 11).Value = 1214.6875
$ws.Cells.Item(113, 13).Value = 955.3125
$ws.Cells.Item(122, 8).Value = 2400.6365
$ws.Cells.Item(122, 10).Value = 3821
$ws.Cells.Item(122, 12).Value = 11463
$ws.Cells.Item(122, 14).Value = -16363
$ws.Cells.Item(132, 8).Value = 2547.6296
$ws.Cells.Item(132, 9).Value = 2204.0952
$ws.Cells.Item(132, 10).Value = 3750
$ws.Cells.Item(132, 11).Value = 6612.285600000001
$ws.Cells.Item(132, 12).Value = 11250
$ws.Cells.Item(132, 13).Value = -4082.285600000001
$ws.Cells.Item(132, 14).Value = -16310
$ws.Cells.Item(134, 8).Value = 6351.1055
$ws.Cells.Item(134, 9).Value = 6000.6665
$ws.Cells.Item(134, 11).Value = 18001.9995
$ws.Cells.Item(134, 13).Value = -15466.9995
$ws.Cells.Item(136, 8).Value = 2538.3125
$ws.Cells.Item(136, 9).Value = 1656.2858
$ws.Cells.Item(136, 10).Value = 4222.1816
$ws.Cells.Item(136, 11).Value = 4968.857400000001
$ws.Cells.Item(136, 12).Value = 12666.5448
$ws.Cells.Item(136, 13).Value = -2418.857400000001
$ws.Cells.Item(136, 14).Value = -17766.5448

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(134, 8).Value = 16668244
$ws.Cells.Item(134, 9).Value = 16668244
$ws.Cells.Item(134, 10).Value = 0
$ws.Cells.Item(134, 11).Value = 50004732
$ws.Cells.Item(134, 12).Value = 0
$ws.Cells.Item(134, 13).Value = -49999662
$ws.Cells.Item(134, 14).ClearContents()

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(132, 8).Value = 4559.727
$ws.Cells.Item(132, 9).Value = 1265.75
$ws.Cells.Item(132, 11).Value = 3797.25
$ws.Cells.Item(132, 13).Value = -1267.25

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 3271
$ws.Cells.Item(7, 9).Value = 3330
$ws.Cells.Item(7, 10).Value = 3049.75
$ws.Cells.Item(7, 11).Value = 3330
$ws.Cells.Item(7, 12).Value = 3049.75
$ws.Cells.Item(7, 13).Value = -3218
$ws.Cells.Item(7, 14).Value = -3273.75
$ws.Cells.Item(68, 8).Value = 2374.2856
$ws.Cells.Item(68, 9).Value = 2037.7778
$ws.Cells.Item(68, 10).Value = 2980
$ws.Cells.Item(68, 11).Value = 2037.7778
$ws.Cells.Item(68, 12).Value = 2980
$ws.Cells.Item(68, 13).Value = -1288.7778
$ws.Cells.Item(68, 14).Value = -4478
$ws.Cells.Item(71, 8).Value = 2374.2856
$ws.Cells.Item(71, 9).Value = 2037.7778
$ws.Cells.Item(71, 10).Value = 2980
$ws.Cells.Item(71, 11).Value = 10188.889
$ws.Cells.Item(71, 12).Value = 14900
$ws.Cells.Item(71, 13).Value = -6444.889000000001
$ws.Cells.Item(71, 14).Value = -22388
$ws.Cells.Item(122, 8).Value = 4415.3335
$ws.Cells.Item(122, 9).Value = 4453.091
$ws.Cells.Item(122, 11).Value = 13359.273
$ws.Cells.Item(122, 13).Value = -10909.273
$ws.Cells.Item(126, 8).Value = 3271
$ws.Cells.Item(126, 9).Value = 3330
$ws.Cells.Item(126, 10).Value = 3049.75
$ws.Cells.Item(126, 11).Value = 9990
$ws.Cells.Item(126, 12).Value = 9149.25
$ws.Cells.Item(126, 13).Value = -7520
$ws.Cells.Item(126, 14).Value = -14089.25
$ws.Cells.Item(132, 8).Value = 5162.4
$ws.Cells.Item(132, 9).Value = 4309.1577
$ws.Cells.Item(132, 10).Value = 6636.1816
$ws.Cells.Item(132, 11).Value = 12927.4731
$ws.Cells.Item(132, 12).Value = 19908.5448
$ws.Cells.Item(132, 13).Value = -10397.4731
$ws.Cells.Item(132, 14).Value = -24968.5448
$ws.Cells.Item(136, 8).Value = 3605.6667
$ws.Cells.Item(136, 9).Value = 3456.5454
$ws.Cells.Item(136, 11).Value = 10369.6362
$ws.Cells.Item(136, 13).Value = -7819.636200000001

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(62, 8).Value = 4839
$ws.Cells.Item(62, 9).Value = 4028.6296
$ws.Cells.Item(62, 11).Value = 4028.6296
$ws.Cells.Item(62, 13).Value = -3404.6296
$ws.Cells.Item(65, 8).Value = 4839
$ws.Cells.Item(65, 9).Value = 4028.6296
$ws.Cells.Item(65, 11).Value = 20143.148
$ws.Cells.Item(65, 13).Value = -17023.148
$ws.Cells.Item(113, 8).Value = 513.16
$ws.Cells.Item(113, 9).Value = 490.625
$ws.Cells.Item(113, 10).Value = 553.2222
$ws.Cells.Item(113, 11).Value = 1471.875
$ws.Cells.Item(113, 12).Value = 1659.6666
$ws.Cells.Item(113, 13).Value = 698.125
$ws.Cells.Item(113, 14).Value = -5999.6666
$ws.Cells.Item(132, 8).Value = 2537.3667
$ws.Cells.Item(132, 9).Value = 2773.32
$ws.Cells.Item(132, 10).Value = 1357.6
$ws.Cells.Item(132, 11).Value = 8319.960000000001
$ws.Cells.Item(132, 12).Value = 4072.8
$ws.Cells.Item(132, 13).Value = -5789.960000000001
$ws.Cells.Item(132, 14).Value = -9132.799999999999
